$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3461, 3719, 3719, 4195, 4400, 4472, 4528, 4958, 5002, 5002, 5021, 5021, 5021, 5217)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
